# Apply the edit described by the commit:
#   - enter "ffdd" into A1 and "kccc" into B1 (as text -> becomes shared strings)
#   - widen columns A and B to fit the new content
#   - leave the selection / active cell on B1

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "ffdd"
$ws.Range("B1").Value = "kccc"

$ws.Columns.Item(1).ColumnWidth = 13.16
$ws.Columns.Item(2).ColumnWidth = 13.3

$ws.Range("B1").Select() | Out-Null
